$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 27 (pushes existing rows 27-71 down to 28-72,
# carrying the date-format style from the old D27 cell along with the shift).
$ws.Rows.Item(27).Insert()

# Populate the new row 27 with the latest weekly price observation.
$ws.Range("A27").Value = 1
$ws.Range("B27").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C27").Value = "Arica y Parinacota"
$ws.Range("D27").Value = 44540
$ws.Range("E27").Value = 15
$ws.Range("F27").Value = "Fruta"
$ws.Range("G27").Value = 100106
$ws.Range("H27").Value = "Oleaginosos"
$ws.Range("I27").Value = 100106002
$ws.Range("J27").Value = "Palta"
$ws.Range("K27").Value = "Hass"
$ws.Range("L27").Value = "Primera"
$ws.Range("M27").Value = 400
$ws.Range("N27").Value = 23000
$ws.Range("O27").Value = 24000
$ws.Range("P27").Value = 23500
$ws.Range("Q27").Value = "$/bandeja 10 kilos"
$ws.Range("R27").Value = "Perú"
$ws.Range("S27").Value = 2350
$ws.Range("T27").Value = 10
